# Applies the "Deploying to gh-pages ... Alvearie/alvearie-fhir-ig@8e4a450c..." edit:
#   - Metadata sheet: Version 5.0.0 -> 6.0.0
#   - Metadata sheet: Date 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
#   - Metadata sheet: Publisher value "" -> "Alvearie Team"
#   - Metadata sheet: duplicate "Contact / No display for ContactDetail" row collapsed,
#     replaced by a new "Jurisdiction / United States of America" row
#   - Elements sheet: root Extension row's Short/Definition ("Extension" / "An Extension")
#     replaced with the StructureDefinition's own Title/Description
#     ("Capitation Arrangement" / "Code for the capitation arrangement")

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# There were two identical "Contact / No display for ContactDetail" rows
# (rows 10 and 11). Remove one of them entirely (row 11) ...
$meta.Rows(11).Delete()

# ... and turn the remaining one (now still row 10) into the new
# "Jurisdiction / United States of America" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------

# Row 2 is the root "Extension" element. Its Short/Definition columns (K/L)
# generically said "Extension" / "An Extension" - replace with the actual
# title/description of this extension.
$elements.Range("K2").Value = "Capitation Arrangement"
$elements.Range("L2").Value = "Code for the capitation arrangement"
